$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("other_data_caseI")

# New headers for columns AD:AI ("... Cost for Failures" metrics)
$headers = @(
    "ATTA H Cost for Failures",
    "ATTA R Cost for Failures",
    "Random H Cost for Failures",
    "Random R Cost for Failures",
    "Tsarouchi MIN H Cost for Failures",
    "Tsarouchi MIN R Cost for Failures"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 30 + $i   # AD = column 30
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# New data rows 2..11 for columns AD:AI
$data = @(
    @(2.0843366441311799, 11.3005853723268, 16.432125115786501, 12.559642086380199, 9.3898869019654594, 3.9959030524759398),
    @(1.33725748000518, 14.3335699508559, 21.722011581827299, 12.889650581358801, 10.2739830714657, 4.7325541979824903),
    @(2.3566511075727701, 13.658935975225599, 19.8765132676115, 13.3454528397969, 11.2566146279178, 4.5002690207741001),
    @(2.0798293450997001, 13.178863475239, 22.815951038125299, 14.209063805568301, 12.204985795034499, 4.9240618668343901),
    @(0.82330166017134998, 13.479453889128299, 23.0478389263879, 11.3044702886891, 12.8660490073537, 4.3606936753822696),
    @(1.93700779548671, 12.4772936183819, 23.020128270590899, 11.630967747175299, 10.431218249683001, 4.2552139952902799),
    @(1.2758997942449699, 14.049310205823399, 19.923239031373701, 13.260277285523101, 13.148432820491999, 4.5391094966868897),
    @(1.35007971880422, 15.092534410311201, 23.6100964275333, 13.4076155817519, 11.987769179998899, 4.9435165177143201),
    @(2.3816939623441198, 12.0065924121413, 18.175618560538499, 11.1517607286002, 9.9595832081362108, 4.5272965123244502),
    @(0.65111910368157, 14.149850474089799, 24.371990196236901, 12.743832954299799, 11.0711078085212, 4.6526274827717504)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $col = 30 + $c   # AD = column 30
        $ws.Cells.Item($excelRow, $col).Value = $rowVals[$c]
    }
}

# Move the active selection onto the newly added data (matches the saved view)
$ws.Range("AA5").Select()
